$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header in H1, matching the style of the existing headers (B1:G1)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Fill in the "Save" column values for rows 2-9
$saveValues = @(0, 0, 0, 0, 1, 1, 0, 0)
for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
